# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the newly scraped counts (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (rows 2-11) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 8567
$ws1.Range("F3").Value = 88
$ws1.Range("F4").Value = 230
$ws1.Range("F5").Value = 86
$ws1.Range("F6").Value = 1289
$ws1.Range("F9").Value = 32
$ws1.Range("F10").Value = 231
$ws1.Range("F11").Value = 66

# --- Sheet "全部类型" (rows 2-12) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 8567
$ws4.Range("F3").Value = 88
$ws4.Range("F4").Value = 230
$ws4.Range("F5").Value = 86
$ws4.Range("F6").Value = 1289
$ws4.Range("F10").Value = 32
$ws4.Range("F11").Value = 231
$ws4.Range("F12").Value = 66

$wb.Save()
